$wb = $excel.ActiveWorkbook

# --- HARDWARE_MANAGEMENT (sheet3): add BMC port rows for the river compute nodes ---
$wsHW = $wb.Worksheets.Item("HARDWARE_MANAGEMENT")

$hwRows = @(
    @{ Row=25; J="cn01"; K="x3002"; L="u15"; M="bmc"; O=1; P="sw-leaf-bmc-001"; Q="x3000"; R="u37"; T=11 },
    @{ Row=26; J="cn02"; K="x3002"; L="u16"; M="bmc"; O=1; P="sw-leaf-bmc-001"; Q="x3000"; R="u37"; T=12 },
    @{ Row=27; J="cn03"; K="x3002"; L="u17"; M="bmc"; O=1; P="sw-leaf-bmc-001"; Q="x3000"; R="u37"; T=13 },
    @{ Row=28; J="cn04"; K="x3002"; L="u18"; M="bmc"; O=1; P="sw-leaf-bmc-001"; Q="x3000"; R="u37"; T=14 }
)

foreach ($r in $hwRows) {
    $wsHW.Range("J" + $r.Row).Value = $r.J
    $wsHW.Range("K" + $r.Row).Value = $r.K
    $wsHW.Range("L" + $r.Row).Value = $r.L
    $wsHW.Range("M" + $r.Row).Value = $r.M
    $wsHW.Range("O" + $r.Row).Value = $r.O
    $wsHW.Range("P" + $r.Row).Value = $r.P
    $wsHW.Range("Q" + $r.Row).Value = $r.Q
    $wsHW.Range("R" + $r.Row).Value = $r.R
    $wsHW.Range("T" + $r.Row).Value = $r.T
}

# --- COMPUTE_NODES (sheet4): add BMC port rows for the river compute nodes ---
$wsCN = $wb.Worksheets.Item("COMPUTE_NODES")

$cnRows = @(
    @{ Row=24; J="cn01"; K="x3002"; L="u15"; O=1; P="sw-leaf-bmc-001"; Q="x3000"; R="u37"; T=24 },
    @{ Row=25; J="cn02"; K="x3002"; L="u16"; O=1; P="sw-leaf-bmc-001"; Q="x3000"; R="u37"; T=25 },
    @{ Row=26; J="cn03"; K="x3002"; L="u17"; O=1; P="sw-leaf-bmc-001"; Q="x3000"; R="u37"; T=26 },
    @{ Row=27; J="cn04"; K="x3002"; L="u18"; O=1; P="sw-leaf-bmc-001"; Q="x3000"; R="u37"; T=27 }
)

foreach ($r in $cnRows) {
    $wsCN.Range("J" + $r.Row).Value = $r.J
    $wsCN.Range("K" + $r.Row).Value = $r.K
    $wsCN.Range("L" + $r.Row).Value = $r.L
    $wsCN.Range("O" + $r.Row).Value = $r.O
    $wsCN.Range("P" + $r.Row).Value = $r.P
    $wsCN.Range("Q" + $r.Row).Value = $r.Q
    $wsCN.Range("R" + $r.Row).Value = $r.R
    $wsCN.Range("T" + $r.Row).Value = $r.T
}

# --- Update selections to match the new data extents ---
$wsHW.Range("J25:T28").Select()

# COMPUTE_NODES becomes the active/selected tab last, matching the saved view state
$wsCN.Range("J24:T27").Select()
$wsCN.Activate()

Write-Host "edit complete"
